$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that sits at the end of the paragraph
#    describing the "100th count" (it is no longer the last-edit location
#    once we add the new text below).
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 2) Find the paragraph that ends in "... on the 1000th count" and append the
#    two new sentences, with a fresh "_GoBack" bookmark sitting between them
#    (this mirrors Word automatically re-stamping the last-edit bookmark at
#    the new insertion point).
# ---------------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Ring Finger on the 1000*count*") {
        $targetPara = $p
        break
    }
}

$endPos = $targetPara.Range.End - 1

$runE = ": I certainly didn" + [char]0x2019 + "t think I could count until 1000 and actually get an accurate result so I decided that I would try to find"
$runF = " a pattern in the counting. I ultimately found that every 50 fingers would land on either the First Finger or the Ring Finger. So I divided 1000 by 50 and got 20. I then varied between the Right and Ring finger and ultimately concluded that the Ring Ringer would be count 1000."

# Insert the first new run right after the existing " count" run.
$r = $d.Range($endPos, $endPos)
$r.InsertAfter($runE)
$r.Font.Size = 14

# Remember the boundary between the two new runs - this is where the
# "_GoBack" bookmark needs to land.
$bmPos = $r.End

# Insert the second new run right after the first one.
$r.Collapse(0)
$r.InsertAfter($runF)
$r.Font.Size = 14

# Drop the new "_GoBack" bookmark between the two freshly inserted runs.
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
